# Auto-generated edit script applying scheduled market-price refresh
# to the Leve profit tables (columns H-N) across all 8 crafting-class sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 33589.25
$ws.Range("J3").Value = 33589.25
$ws.Range("L3").Value = 33589.25
$ws.Range("N3").Value = -33817.25
$ws.Range("H32").Value = 2199.1
$ws.Range("I32").Value = 838.4
$ws.Range("K32").Value = 838.4
$ws.Range("M32").Value = -512.4
$ws.Range("H102").Value = 33589.25
$ws.Range("J102").Value = 33589.25
$ws.Range("L102").Value = 33589.25
$ws.Range("N102").Value = -40079.25
$ws.Range("H103").Value = 5887.0586
$ws.Range("I103").Value = 727.7778
$ws.Range("J103").Value = 11691.25
$ws.Range("K103").Value = 2183.3334
$ws.Range("L103").Value = 35073.75
$ws.Range("M103").Value = -1597.3334
$ws.Range("N103").Value = -36245.75
$ws.Range("H111").Value = 624.5
$ws.Range("I111").Value = 624.5
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 1873.5
$ws.Range("L111").Value = 0
$ws.Range("M111").ClearContents()
$ws.Range("N111").Value = 1193.5
$ws.Range("H112").Value = 1321.9412
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 1321.9412
$ws.Range("K112").Value = 0
$ws.Range("L112").ClearContents()
$ws.Range("M112").Value = 3965.8236
$ws.Range("N112").Value = -6181.8236
$ws.Range("H132").Value = 307602.38
$ws.Range("I132").Value = 205346.42
$ws.Range("J132").Value = 1002942.8
$ws.Range("K132").Value = 616039.26
$ws.Range("L132").Value = 3008828.4
$ws.Range("M132").Value = -613509.26
$ws.Range("N132").Value = -3013888.4
$ws.Range("H135").Value = 685.625
$ws.Range("I135").Value = 503.6
$ws.Range("J135").Value = 989
$ws.Range("K135").Value = 4532.400000000001
$ws.Range("L135").Value = 8901
$ws.Range("M135").Value = -1997.400000000001
$ws.Range("N135").Value = -13971

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2378.12
$ws.Range("I32").Value = 2091.886
$ws.Range("J32").Value = 3454.9048
$ws.Range("K32").Value = 2091.886
$ws.Range("L32").Value = 3454.9048
$ws.Range("M32").Value = -1804.886
$ws.Range("N32").Value = -4028.9048
$ws.Range("H61").Value = 2036.1177
$ws.Range("I61").Value = 2195.889
$ws.Range("K61").Value = 2195.889
$ws.Range("M61").Value = -1983.889
$ws.Range("H74").Value = 3084.3333
$ws.Range("I74").Value = 2494.5
$ws.Range("J74").Value = 4264
$ws.Range("K74").Value = 2494.5
$ws.Range("L74").Value = 4264
$ws.Range("M74").Value = -1620.5
$ws.Range("N74").Value = -6012
$ws.Range("H77").Value = 3084.3333
$ws.Range("I77").Value = 2494.5
$ws.Range("J77").Value = 4264
$ws.Range("K77").Value = 12472.5
$ws.Range("L77").Value = 21320
$ws.Range("M77").Value = -8104.5
$ws.Range("N77").Value = -30056
$ws.Range("H136").Value = 2036.1177
$ws.Range("I136").Value = 2195.889
$ws.Range("K136").Value = 6587.667
$ws.Range("M136").Value = -4037.667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 33895
$ws.Range("J103").Value = 33895
$ws.Range("L103").Value = 33895
$ws.Range("N103").Value = -36239
$ws.Range("H107").Value = 1962.5
$ws.Range("I107").Value = 1940
$ws.Range("K107").Value = 1940
$ws.Range("M107").Value = -20

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 11670.333
$ws.Range("I32").Value = 12500
$ws.Range("J32").Value = 10011
$ws.Range("K32").Value = 12500
$ws.Range("L32").Value = 10011
$ws.Range("M32").Value = -12184
$ws.Range("N32").Value = -10643

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 16054.333
$ws.Range("I26").Value = 20248.334
$ws.Range("J26").Value = 7666.3335
$ws.Range("K26").Value = 60745.00199999999
$ws.Range("L26").Value = 22999.0005
$ws.Range("M26").Value = -60457.00199999999
$ws.Range("N26").Value = -23575.0005
$ws.Range("H87").Value = 7910.5
$ws.Range("I87").Value = 6852.6
$ws.Range("J87").Value = 13200
$ws.Range("K87").Value = 20557.8
$ws.Range("L87").Value = 39600
$ws.Range("M87").Value = -19309.8
$ws.Range("N87").Value = -42096
$ws.Range("H90").Value = 7910.5
$ws.Range("I90").Value = 6852.6
$ws.Range("J90").Value = 13200
$ws.Range("K90").Value = 61673.4
$ws.Range("L90").Value = 118800
$ws.Range("M90").Value = -55433.4
$ws.Range("N90").Value = -131280
$ws.Range("H107").Value = 6263728
$ws.Range("I107").Value = 522.5161000000001
$ws.Range("J107").Value = 10226164
$ws.Range("K107").Value = 1567.5483
$ws.Range("L107").Value = 30678492
$ws.Range("M107").Value = 352.4516999999998
$ws.Range("N107").Value = -30682332
$ws.Range("H112").Value = 4326.6665
$ws.Range("J112").Value = 4326.6665
$ws.Range("L112").Value = 12979.9995
$ws.Range("N112").Value = -15195.9995
$ws.Range("H114").Value = 66668444
$ws.Range("I114").Value = 200000500
$ws.Range("J114").Value = 2420
$ws.Range("K114").Value = 600001500
$ws.Range("L114").Value = 7260
$ws.Range("M114").Value = -599998246
$ws.Range("N114").Value = -13768
$ws.Range("H131").Value = 6025075.5
$ws.Range("I131").Value = 41667080
$ws.Range("J131").Value = 1074.8169
$ws.Range("K131").Value = 125001240
$ws.Range("L131").Value = 3224.4507
$ws.Range("M131").Value = -124996200
$ws.Range("N131").Value = -13304.4507

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4220.6665
$ws.Range("I132").Value = 2322.111
$ws.Range("J132").Value = 5359.8
$ws.Range("K132").Value = 6966.333
$ws.Range("L132").Value = 16079.4
$ws.Range("M132").Value = -4436.333
$ws.Range("N132").Value = -21139.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3145.1035
$ws.Range("I7").Value = 1406.4375
$ws.Range("K7").Value = 1406.4375
$ws.Range("M7").Value = -1294.4375
$ws.Range("H22").Value = 15688919
$ws.Range("I22").Value = 20917542
$ws.Range("J22").Value = 3050
$ws.Range("K22").Value = 20917542
$ws.Range("L22").Value = 3050
$ws.Range("M22").Value = -20917247
$ws.Range("N22").Value = -3640
$ws.Range("H27").Value = 15688919
$ws.Range("I27").Value = 20917542
$ws.Range("J27").Value = 3050
$ws.Range("K27").Value = 20917542
$ws.Range("L27").Value = 3050
$ws.Range("M27").Value = -20917435
$ws.Range("N27").Value = -3264
$ws.Range("H33").Value = 19929.666
$ws.Range("I33").Value = 17400
$ws.Range("J33").Value = 24989
$ws.Range("K33").Value = 17400
$ws.Range("L33").Value = 24989
$ws.Range("M33").Value = -17110
$ws.Range("N33").Value = -25569
$ws.Range("H61").Value = 1270
$ws.Range("I61").Value = 1238.0667
$ws.Range("J61").Value = 1349.8334
$ws.Range("K61").Value = 1238.0667
$ws.Range("L61").Value = 1349.8334
$ws.Range("M61").Value = -1036.0667
$ws.Range("N61").Value = -1753.8334
$ws.Range("H110").Value = 31333.334
$ws.Range("J110").Value = 31333.334
$ws.Range("L110").Value = 31333.334
$ws.Range("N110").Value = -39513.334
$ws.Range("H113").Value = 1270
$ws.Range("I113").Value = 1238.0667
$ws.Range("J113").Value = 1349.8334
$ws.Range("K113").Value = 1238.0667
$ws.Range("L113").Value = 1349.8334
$ws.Range("M113").Value = 931.9332999999999
$ws.Range("N113").Value = -5689.8334
$ws.Range("H126").Value = 3145.1035
$ws.Range("I126").Value = 1406.4375
$ws.Range("K126").Value = 4219.3125
$ws.Range("M126").Value = -1749.3125
$ws.Range("H136").Value = 3611.0264
$ws.Range("I136").Value = 1575.5834
$ws.Range("J136").Value = 7100.357
$ws.Range("K136").Value = 4726.7502
$ws.Range("L136").Value = 21301.071
$ws.Range("M136").Value = -2176.7502
$ws.Range("N136").Value = -26401.071

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 303.2857
$ws.Range("I113").Value = 320.93332
$ws.Range("J113").Value = 259.16666
$ws.Range("K113").Value = 962.7999599999999
$ws.Range("L113").Value = 777.4999799999999
$ws.Range("M113").Value = 1207.20004
$ws.Range("N113").Value = -5117.49998
$ws.Range("H115").Value = 35047.617
$ws.Range("J115").Value = 35047.617
$ws.Range("L115").Value = 35047.617
$ws.Range("N115").Value = -38181.617
$ws.Range("H132").Value = 8774065
$ws.Range("I132").Value = 1486.1666
$ws.Range("J132").Value = 23812770
$ws.Range("K132").Value = 4458.4998
$ws.Range("L132").Value = 71438310
$ws.Range("M132").Value = -1928.4998
$ws.Range("N132").Value = -71443370
$ws.Range("H136").Value = 4637.4346
$ws.Range("I136").Value = 2164.5833
$ws.Range("J136").Value = 7335.091
$ws.Range("K136").Value = 6493.749899999999
$ws.Range("L136").Value = 22005.273
$ws.Range("M136").Value = -3943.749899999999
$ws.Range("N136").Value = -27105.273

